$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 1697
$ws.Cells.Item(4, 9).Value = 1733.625
$ws.Cells.Item(4, 10).Value = 1623.75
$ws.Cells.Item(4, 11).Value = 1733.625
$ws.Cells.Item(4, 12).Value = 1623.75
$ws.Cells.Item(4, 13).Value = -1619.625
$ws.Cells.Item(4, 14).Value = -1851.75
$ws.Cells.Item(5, 8).Value = 285.33334
$ws.Cells.Item(5, 10).Value = 510
$ws.Cells.Item(5, 12).Value = 510
$ws.Cells.Item(5, 14).Value = -740
$ws.Cells.Item(19, 8).Value = 1940.9286
$ws.Cells.Item(19, 9).Value = 2380
$ws.Cells.Item(19, 11).Value = 2380
$ws.Cells.Item(19, 13).Value = -2205
$ws.Cells.Item(32, 8).Value = 4142.25
$ws.Cells.Item(32, 10).Value = 5229.8335
$ws.Cells.Item(32, 12).Value = 5229.8335
$ws.Cells.Item(32, 14).Value = -5881.8335
$ws.Cells.Item(40, 8).Value = 1999.5
$ws.Cells.Item(40, 9).Value = 1628.2858
$ws.Cells.Item(40, 10).Value = 2288.2222
$ws.Cells.Item(40, 11).Value = 1628.2858
$ws.Cells.Item(40, 12).Value = 2288.2222
$ws.Cells.Item(40, 13).Value = -1453.2858
$ws.Cells.Item(40, 14).Value = -2638.2222
$ws.Cells.Item(51, 8).Value = 3650
$ws.Cells.Item(51, 9).Value = 3500
$ws.Cells.Item(51, 11).Value = 3500
$ws.Cells.Item(51, 13).Value = -3016
$ws.Cells.Item(101, 8).Value = 14286255
$ws.Cells.Item(101, 9).Value = 33333580
$ws.Cells.Item(101, 11).Value = 100000740
$ws.Cells.Item(101, 13).Value = -99999118
$ws.Cells.Item(137, 8).Value = 1624.75
$ws.Cells.Item(137, 9).Value = 1250
$ws.Cells.Item(137, 11).Value = 3750
$ws.Cells.Item(137, 13).Value = -1200
$ws.Cells.Item(138, 8).Value = 5041.25
$ws.Cells.Item(138, 9).Value = 2431.6667
$ws.Cells.Item(138, 10).Value = 5911.1113
$ws.Cells.Item(138, 11).Value = 7295.000100000001
$ws.Cells.Item(138, 12).Value = 17733.3339
$ws.Cells.Item(138, 13).Value = -2155.000100000001
$ws.Cells.Item(138, 14).Value = -28013.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1802.4615
$ws.Cells.Item(61, 9).Value = 1802.4615
$ws.Cells.Item(61, 11).Value = 1802.4615
$ws.Cells.Item(61, 13).Value = -1590.4615
$ws.Cells.Item(63, 8).Value = 3660
$ws.Cells.Item(63, 9).Value = 4153.636
$ws.Cells.Item(63, 11).Value = 4153.636
$ws.Cells.Item(63, 13).Value = -3467.636
$ws.Cells.Item(66, 8).Value = 3660
$ws.Cells.Item(66, 9).Value = 4153.636
$ws.Cells.Item(66, 11).Value = 20768.18
$ws.Cells.Item(66, 13).Value = -17336.18
$ws.Cells.Item(74, 8).Value = 4045
$ws.Cells.Item(74, 9).Value = 3425.625
$ws.Cells.Item(74, 11).Value = 3425.625
$ws.Cells.Item(74, 13).Value = -2551.625
$ws.Cells.Item(77, 8).Value = 4045
$ws.Cells.Item(77, 9).Value = 3425.625
$ws.Cells.Item(77, 11).Value = 17128.125
$ws.Cells.Item(77, 13).Value = -12760.125
$ws.Cells.Item(122, 8).Value = 1401.381
$ws.Cells.Item(122, 9).Value = 1079.4445
$ws.Cells.Item(122, 11).Value = 3238.3335
$ws.Cells.Item(122, 13).Value = -788.3335000000002
$ws.Cells.Item(136, 8).Value = 1802.4615
$ws.Cells.Item(136, 9).Value = 1802.4615
$ws.Cells.Item(136, 11).Value = 5407.3845
$ws.Cells.Item(136, 13).Value = -2857.3845

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 713.1667
$ws.Cells.Item(64, 10).Value = 597.25
$ws.Cells.Item(64, 12).Value = 597.25
$ws.Cells.Item(64, 14).Value = -1047.25
$ws.Cells.Item(67, 8).Value = 713.1667
$ws.Cells.Item(67, 10).Value = 597.25
$ws.Cells.Item(67, 12).Value = 597.25
$ws.Cells.Item(67, 14).Value = -2157.25
$ws.Cells.Item(80, 8).Value = 587.55554
$ws.Cells.Item(80, 9).Value = 660.2
$ws.Cells.Item(80, 10).Value = 496.75
$ws.Cells.Item(80, 11).Value = 660.2
$ws.Cells.Item(80, 12).Value = 496.75
$ws.Cells.Item(80, 13).Value = 337.8
$ws.Cells.Item(80, 14).Value = -2492.75
$ws.Cells.Item(83, 8).Value = 587.55554
$ws.Cells.Item(83, 9).Value = 660.2
$ws.Cells.Item(83, 10).Value = 496.75
$ws.Cells.Item(83, 11).Value = 3301
$ws.Cells.Item(83, 12).Value = 2483.75
$ws.Cells.Item(83, 13).Value = 1691
$ws.Cells.Item(83, 14).Value = -12467.75
$ws.Cells.Item(95, 8).Value = 21265
$ws.Cells.Item(95, 10).Value = 21265
$ws.Cells.Item(95, 12).Value = 21265
$ws.Cells.Item(95, 14).Value = -26757
$ws.Cells.Item(107, 8).Value = 495.8889
$ws.Cells.Item(107, 9).Value = 471.66666
$ws.Cells.Item(107, 10).Value = 544.3333
$ws.Cells.Item(107, 11).Value = 471.66666
$ws.Cells.Item(107, 12).Value = 544.3333
$ws.Cells.Item(107, 13).Value = 1448.33334
$ws.Cells.Item(107, 14).Value = -4384.3333
$ws.Cells.Item(134, 8).Value = 2049
$ws.Cells.Item(134, 9).Value = 1984.9445
$ws.Cells.Item(134, 11).Value = 5954.833500000001
$ws.Cells.Item(134, 13).Value = -3419.833500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2741
$ws.Cells.Item(31, 9).Value = 1663.6666
$ws.Cells.Item(31, 11).Value = 1663.6666
$ws.Cells.Item(31, 13).Value = -1368.6666
$ws.Cells.Item(34, 8).Value = 2741
$ws.Cells.Item(34, 9).Value = 1663.6666
$ws.Cells.Item(34, 11).Value = 1663.6666
$ws.Cells.Item(34, 13).Value = -1461.6666
$ws.Cells.Item(58, 8).Value = 1947.3846
$ws.Cells.Item(58, 9).Value = 1739.3334
$ws.Cells.Item(58, 11).Value = 1739.3334
$ws.Cells.Item(58, 13).Value = -1536.3334
$ws.Cells.Item(94, 8).Value = 65605.164
$ws.Cells.Item(94, 9).Value = 125128.555
$ws.Cells.Item(94, 10).Value = 6081.778
$ws.Cells.Item(94, 11).Value = 125128.555
$ws.Cells.Item(94, 12).Value = 6081.778
$ws.Cells.Item(94, 13).Value = -124677.555
$ws.Cells.Item(94, 14).Value = -6983.778
$ws.Cells.Item(106, 8).Value = 27750
$ws.Cells.Item(106, 10).Value = 27750
$ws.Cells.Item(106, 12).Value = 27750
$ws.Cells.Item(106, 14).Value = -30274
$ws.Cells.Item(122, 8).Value = 1908
$ws.Cells.Item(122, 9).Value = 1830
$ws.Cells.Item(122, 11).Value = 5490
$ws.Cells.Item(122, 13).Value = -3040
$ws.Cells.Item(132, 8).Value = 3408
$ws.Cells.Item(132, 9).Value = 3481.3572
$ws.Cells.Item(132, 10).Value = 3065.6667
$ws.Cells.Item(132, 11).Value = 10444.0716
$ws.Cells.Item(132, 12).Value = 9197.000100000001
$ws.Cells.Item(132, 13).Value = -7914.071599999999
$ws.Cells.Item(132, 14).Value = -14257.0001
$ws.Cells.Item(134, 8).Value = 2376.8845
$ws.Cells.Item(134, 9).Value = 2332.8262
$ws.Cells.Item(134, 11).Value = 6998.4786
$ws.Cells.Item(134, 13).Value = -4463.4786
$ws.Cells.Item(136, 8).Value = 1947.3846
$ws.Cells.Item(136, 9).Value = 1739.3334
$ws.Cells.Item(136, 11).Value = 5218.0002
$ws.Cells.Item(136, 13).Value = -2668.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 1252842.1
$ws.Cells.Item(129, 9).Value = 2443.5
$ws.Cells.Item(129, 10).Value = 2503240.8
$ws.Cells.Item(129, 11).Value = 7330.5
$ws.Cells.Item(129, 12).Value = 7509722.399999999
$ws.Cells.Item(129, 13).Value = -2330.5
$ws.Cells.Item(129, 14).Value = -7519722.399999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 3318.125
$ws.Cells.Item(107, 9).Value = 1497.5
$ws.Cells.Item(107, 11).Value = 1497.5
$ws.Cells.Item(107, 13).Value = 422.5
$ws.Cells.Item(122, 8).Value = 1028.75
$ws.Cells.Item(122, 9).Value = 1028.75
$ws.Cells.Item(122, 11).Value = 3086.25
$ws.Cells.Item(122, 13).Value = -636.25
$ws.Cells.Item(126, 8).Value = 2399.1428
$ws.Cells.Item(126, 9).Value = 2299.5
$ws.Cells.Item(126, 11).Value = 6898.5
$ws.Cells.Item(126, 13).Value = -4428.5
$ws.Cells.Item(132, 8).Value = 2479.875
$ws.Cells.Item(132, 9).Value = 2479.875
$ws.Cells.Item(132, 11).Value = 7439.625
$ws.Cells.Item(132, 13).Value = -4909.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 10799.2
$ws.Cells.Item(7, 9).Value = 11999
$ws.Cells.Item(7, 10).Value = 6000
$ws.Cells.Item(7, 11).Value = 11999
$ws.Cells.Item(7, 12).Value = 6000
$ws.Cells.Item(7, 13).Value = -11887
$ws.Cells.Item(7, 14).Value = -6224
$ws.Cells.Item(97, 8).Value = 18000
$ws.Cells.Item(97, 10).Value = 18000
$ws.Cells.Item(97, 12).Value = 18000
$ws.Cells.Item(97, 14).Value = -19982
$ws.Cells.Item(126, 8).Value = 10799.2
$ws.Cells.Item(126, 9).Value = 11999
$ws.Cells.Item(126, 10).Value = 6000
$ws.Cells.Item(126, 11).Value = 35997
$ws.Cells.Item(126, 12).Value = 18000
$ws.Cells.Item(126, 13).Value = -33527
$ws.Cells.Item(126, 14).Value = -22940
$ws.Cells.Item(132, 8).Value = 8864.666999999999
$ws.Cells.Item(132, 9).Value = 11447.167
$ws.Cells.Item(132, 11).Value = 34341.501
$ws.Cells.Item(132, 13).Value = -31811.501

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(6, 8).Value = 3980.2
$ws.Cells.Item(6, 9).Value = 3500.3333
$ws.Cells.Item(6, 10).Value = 4700
$ws.Cells.Item(6, 11).Value = 3500.3333
$ws.Cells.Item(6, 12).Value = 4700
$ws.Cells.Item(6, 13).Value = -3385.3333
$ws.Cells.Item(6, 14).Value = -4930
$ws.Cells.Item(8, 8).Value = 457027
$ws.Cells.Item(8, 9).Value = 570033.75
$ws.Cells.Item(8, 10).Value = 5000
$ws.Cells.Item(8, 11).Value = 570033.75
$ws.Cells.Item(8, 12).Value = 5000
$ws.Cells.Item(8, 13).Value = -569893.75
$ws.Cells.Item(8, 14).Value = -5280
$ws.Cells.Item(11, 8).Value = 128500
$ws.Cells.Item(11, 10).Value = 2000
$ws.Cells.Item(11, 12).Value = 2000
$ws.Cells.Item(11, 14).Value = -2284
$ws.Cells.Item(42, 8).Value = 14999.5
$ws.Cells.Item(42, 9).Value = 14999
$ws.Cells.Item(42, 11).Value = 14999
$ws.Cells.Item(42, 13).Value = -14621
$ws.Cells.Item(81, 8).Value = 1113608
$ws.Cells.Item(81, 10).Value = 1669830.5
$ws.Cells.Item(81, 12).Value = 3339661
$ws.Cells.Item(81, 14).Value = -3341783
$ws.Cells.Item(84, 8).Value = 1113608
$ws.Cells.Item(84, 10).Value = 1669830.5
$ws.Cells.Item(84, 12).Value = 16698305
$ws.Cells.Item(84, 14).Value = -16708913
$ws.Cells.Item(107, 8).Value = 377.57144
$ws.Cells.Item(107, 9).Value = 270.9091
$ws.Cells.Item(107, 10).Value = 768.6667
$ws.Cells.Item(107, 11).Value = 812.7273
$ws.Cells.Item(107, 12).Value = 2306.0001
$ws.Cells.Item(107, 13).Value = 1107.2727
$ws.Cells.Item(107, 14).Value = -6146.0001
$ws.Cells.Item(132, 8).Value = 1937.4117
$ws.Cells.Item(132, 9).Value = 1975.7333
$ws.Cells.Item(132, 11).Value = 5927.199900000001
$ws.Cells.Item(132, 13).Value = -3397.199900000001

Write-Output "Applied all changes"